$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1014.1667
$ws.Range("I6").Value = 1014.1667
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 3042.5001
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -2930.5001
$ws.Range("N6").ClearContents()

$ws.Range("H51").Value = 1399.9
$ws.Range("I51").Value = 999.5
$ws.Range("K51").Value = 999.5
$ws.Range("M51").Value = -515.5

$ws.Range("H131").Value = 9840
$ws.Range("J131").Value = 12000
$ws.Range("L131").Value = 36000
$ws.Range("N131").Value = -46080

$ws.Range("H132").Value = 1233.5
$ws.Range("I132").Value = 1233.5
$ws.Range("K132").Value = 3700.5
$ws.Range("M132").Value = -1170.5

$ws.Range("H137").Value = 3944.9333
$ws.Range("I137").Value = 3570.4546
$ws.Range("K137").Value = 10711.3638
$ws.Range("M137").Value = -8161.363799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1210.6666
$ws.Range("I61").Value = 1299.4286
$ws.Range("K61").Value = 1299.4286
$ws.Range("M61").Value = -1087.4286

$ws.Range("H63").Value = 10165.833
$ws.Range("J63").Value = 14998.75
$ws.Range("L63").Value = 14998.75
$ws.Range("N63").Value = -16370.75

$ws.Range("H66").Value = 10165.833
$ws.Range("J66").Value = 14998.75
$ws.Range("L66").Value = 74993.75
$ws.Range("N66").Value = -81857.75

$ws.Range("H119").Value = 100000
$ws.Range("J119").Value = 100000
$ws.Range("L119").Value = 100000
$ws.Range("N119").Value = -109676

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 4209.5
$ws.Range("I132").Value = 3766.3333
$ws.Range("J132").Value = 4874.25
$ws.Range("K132").Value = 11298.9999
$ws.Range("L132").Value = 14622.75
$ws.Range("M132").Value = -8768.999899999999
$ws.Range("N132").Value = -19682.75

$ws.Range("H136").Value = 1210.6666
$ws.Range("I136").Value = 1299.4286
$ws.Range("K136").Value = 3898.2858
$ws.Range("M136").Value = -1348.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9279.4
$ws.Range("J20").Value = 13999.333
$ws.Range("L20").Value = 13999.333
$ws.Range("N20").Value = -14493.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4655.25
$ws.Range("I58").Value = 4655.25
$ws.Range("K58").Value = 4655.25
$ws.Range("M58").Value = -4452.25

$ws.Range("H132").Value = 1285.7273
$ws.Range("I132").Value = 783.4706
$ws.Range("K132").Value = 2350.4118
$ws.Range("M132").Value = 179.5882000000001

$ws.Range("H134").Value = 3489.4783
$ws.Range("I134").Value = 3031.8333
$ws.Range("K134").Value = 9095.499899999999
$ws.Range("M134").Value = -6560.499899999999

$ws.Range("H136").Value = 4655.25
$ws.Range("I136").Value = 4655.25
$ws.Range("K136").Value = 13965.75
$ws.Range("M136").Value = -11415.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4344.2856
$ws.Range("J4").Value = 70
$ws.Range("L4").Value = 210
$ws.Range("N4").Value = -434

$ws.Range("H10").Value = 1242.8
$ws.Range("J10").Value = 2499
$ws.Range("L10").Value = 7497
$ws.Range("N10").Value = -7775

$ws.Range("H32").Value = 1025
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1025
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 3075
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3641

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 769.4286
$ws.Range("I107").Value = 1365.6666
$ws.Range("J107").Value = 322.25
$ws.Range("K107").Value = 1365.6666
$ws.Range("L107").Value = 322.25
$ws.Range("M107").Value = 554.3334
$ws.Range("N107").Value = -4162.25

$ws.Range("H126").Value = 4658.3335
$ws.Range("I126").Value = 4658.3335
$ws.Range("K126").Value = 13975.0005
$ws.Range("M126").Value = -11505.0005

$ws.Range("H132").Value = 1481.25
$ws.Range("I132").Value = 1090.3889
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 3271.1667
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -741.1666999999998
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5040.7144
$ws.Range("I22").Value = 3408.077
$ws.Range("J22").Value = 7693.75
$ws.Range("K22").Value = 3408.077
$ws.Range("L22").Value = 7693.75
$ws.Range("M22").Value = -3113.077
$ws.Range("N22").Value = -8283.75

$ws.Range("H27").Value = 5040.7144
$ws.Range("I27").Value = 3408.077
$ws.Range("J27").Value = 7693.75
$ws.Range("K27").Value = 3408.077
$ws.Range("L27").Value = 7693.75
$ws.Range("M27").Value = -3301.077
$ws.Range("N27").Value = -7907.75

$ws.Range("H76").Value = 87439
$ws.Range("J76").Value = 87439
$ws.Range("L76").Value = 87439
$ws.Range("N76").Value = -88115

$ws.Range("H79").Value = 87439
$ws.Range("J79").Value = 87439
$ws.Range("L79").Value = 87439
$ws.Range("N79").Value = -89779

$ws.Range("H106").Value = 38232.832
$ws.Range("J106").Value = 38232.832
$ws.Range("L106").Value = 38232.832
$ws.Range("N106").Value = -40756.832

$ws.Range("H132").Value = 3769.125
$ws.Range("J132").Value = 4256.4287
$ws.Range("L132").Value = 12769.2861
$ws.Range("N132").Value = -17829.2861

$ws.Range("H136").Value = 3407.0557
$ws.Range("I136").Value = 3308.1333
$ws.Range("K136").Value = 9924.3999
$ws.Range("M136").Value = -7374.3999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

$ws.Range("H62").Value = 6705.5557
$ws.Range("I62").Value = 2833.3333
$ws.Range("J62").Value = 8641.666999999999
$ws.Range("K62").Value = 2833.3333
$ws.Range("L62").Value = 8641.666999999999
$ws.Range("M62").Value = -2209.3333
$ws.Range("N62").Value = -9889.666999999999

$ws.Range("H65").Value = 6705.5557
$ws.Range("I65").Value = 2833.3333
$ws.Range("J65").Value = 8641.666999999999
$ws.Range("K65").Value = 14166.6665
$ws.Range("L65").Value = 43208.335
$ws.Range("M65").Value = -11046.6665
$ws.Range("N65").Value = -49448.335

$ws.Range("H132").Value = 2046.5
$ws.Range("J132").Value = 2995.75
$ws.Range("L132").Value = 8987.25
$ws.Range("N132").Value = -14047.25

$ws.Range("H136").Value = 1658.1111
$ws.Range("I136").Value = 1658.1111
$ws.Range("K136").Value = 4974.3333
$ws.Range("M136").Value = -2424.3333
